# Applies cached-value updates to the Leve-profit tracking tables
# (columns H-N) across all eight crafting-job sheets, per the
# scheduled price-refresh run.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 92
$ws.Range("H92").Value = 1894.4706
$ws.Range("I92").Value = 1680.4
$ws.Range("J92").Value = 3500
$ws.Range("K92").Value = 1680.4
$ws.Range("L92").Value = 3500
$ws.Range("M92").Value = -432.4000000000001
$ws.Range("N92").Value = -5996
# Row 98
$ws.Range("H98").Value = 2536.1482
$ws.Range("I98").Value = 1394.8334
$ws.Range("J98").Value = 11666.667
$ws.Range("K98").Value = 1394.8334
$ws.Range("L98").Value = 11666.667
$ws.Range("M98").Value = 103.1666
$ws.Range("N98").Value = -14662.667
# Row 101
$ws.Range("H101").Value = 666.6667
$ws.Range("I101").Value = 284
$ws.Range("J101").Value = 858
$ws.Range("K101").Value = 852
$ws.Range("L101").Value = 2574
$ws.Range("M101").Value = 770
$ws.Range("N101").Value = -5818
# Row 116
$ws.Range("H116").Value = 352670.44
$ws.Range("J116").Value = 10011.85
$ws.Range("L116").Value = 10011.85
$ws.Range("N116").Value = -16895.85
# Row 122
$ws.Range("H122").Value = 2536.1482
$ws.Range("I122").Value = 1394.8334
$ws.Range("J122").Value = 11666.667
$ws.Range("K122").Value = 4184.5002
$ws.Range("L122").Value = 35000.001
$ws.Range("M122").Value = -1734.5002
$ws.Range("N122").Value = -39900.001
# Row 129
$ws.Range("H129").Value = 573.8889
$ws.Range("I129").Value = 390.14285
$ws.Range("J129").Value = 1217
$ws.Range("K129").Value = 1170.42855
$ws.Range("L129").Value = 3651
$ws.Range("M129").Value = 3829.57145
$ws.Range("N129").Value = -13651
# Row 132
$ws.Range("H132").Value = 62512376
$ws.Range("I132").Value = 100015300
$ws.Range("J132").Value = 7502
$ws.Range("K132").Value = 300045900
$ws.Range("L132").Value = 22506
$ws.Range("M132").Value = -300043370
$ws.Range("N132").Value = -27566
# Row 137
$ws.Range("H137").Value = 3696.653
$ws.Range("I137").Value = 3594.25
$ws.Range("J137").Value = 3794.96
$ws.Range("K137").Value = 10782.75
$ws.Range("L137").Value = 11384.88
$ws.Range("M137").Value = -8232.75
$ws.Range("N137").Value = -16484.88
# Row 138
$ws.Range("H138").Value = 3061.899
$ws.Range("I138").Value = 458.1613
$ws.Range("J138").Value = 4248.897
$ws.Range("K138").Value = 1374.4839
$ws.Range("L138").Value = 12746.691
$ws.Range("M138").Value = 3765.5161
$ws.Range("N138").Value = -23026.691
# Row 141
$ws.Range("H141").Value = 4058.5833
$ws.Range("I141").Value = 4144.2354
$ws.Range("J141").Value = 2602.5
$ws.Range("K141").Value = 12432.7062
$ws.Range("L141").Value = 7807.5
$ws.Range("M141").Value = -7252.706199999999
$ws.Range("N141").Value = -18167.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4074.7036
$ws.Range("I32").Value = 4188.5615
$ws.Range("K32").Value = 4188.5615
$ws.Range("M32").Value = -3901.5615
# Row 61
$ws.Range("H61").Value = 2038.4
$ws.Range("I61").Value = 1538.8889
$ws.Range("J61").Value = 2447.0908
$ws.Range("K61").Value = 1538.8889
$ws.Range("L61").Value = 2447.0908
$ws.Range("M61").Value = -1326.8889
$ws.Range("N61").Value = -2871.0908
# Row 74
$ws.Range("H74").Value = 3258.257
$ws.Range("I74").Value = 3314.6785
$ws.Range("J74").Value = 3032.5715
$ws.Range("K74").Value = 3314.6785
$ws.Range("L74").Value = 3032.5715
$ws.Range("M74").Value = -2440.6785
$ws.Range("N74").Value = -4780.5715
# Row 77
$ws.Range("H77").Value = 3258.257
$ws.Range("I77").Value = 3314.6785
$ws.Range("J77").Value = 3032.5715
$ws.Range("K77").Value = 16573.3925
$ws.Range("L77").Value = 15162.8575
$ws.Range("M77").Value = -12205.3925
$ws.Range("N77").Value = -23898.8575
# Row 132
$ws.Range("H132").Value = 4233.25
$ws.Range("I132").Value = 2619.2856
$ws.Range("J132").Value = 5488.5557
$ws.Range("K132").Value = 7857.8568
$ws.Range("L132").Value = 16465.6671
$ws.Range("M132").Value = -5327.8568
$ws.Range("N132").Value = -21525.6671
# Row 136
$ws.Range("H136").Value = 2038.4
$ws.Range("I136").Value = 1538.8889
$ws.Range("J136").Value = 2447.0908
$ws.Range("K136").Value = 4616.6667
$ws.Range("L136").Value = 7341.2724
$ws.Range("M136").Value = -2066.6667
$ws.Range("N136").Value = -12441.2724

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 140
$ws.Range("H140").Value = 65762.5
$ws.Range("J140").Value = 65762.5
$ws.Range("L140").Value = 65762.5
$ws.Range("N140").Value = -76122.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6581148
$ws.Range("I31").Value = 1291.1538
$ws.Range("J31").Value = 20837504
$ws.Range("K31").Value = 1291.1538
$ws.Range("L31").Value = 20837504
$ws.Range("M31").Value = -996.1538
$ws.Range("N31").Value = -20838094
# Row 34
$ws.Range("H34").Value = 6581148
$ws.Range("I34").Value = 1291.1538
$ws.Range("J34").Value = 20837504
$ws.Range("K34").Value = 1291.1538
$ws.Range("L34").Value = 20837504
$ws.Range("M34").Value = -1089.1538
$ws.Range("N34").Value = -20837908
# Row 58
$ws.Range("H58").Value = 1641.8438
$ws.Range("I58").Value = 1504.473
$ws.Range("K58").Value = 1504.473
$ws.Range("M58").Value = -1301.473
# Row 68
$ws.Range("H68").Value = 52159.625
$ws.Range("J68").Value = 52159.625
$ws.Range("L68").Value = 52159.625
$ws.Range("N68").Value = -53657.625
# Row 71
$ws.Range("H71").Value = 52159.625
$ws.Range("J71").Value = 52159.625
$ws.Range("L71").Value = 156478.875
$ws.Range("N71").Value = -163966.875
# Row 116
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
# Row 132
$ws.Range("H132").Value = 2969.125
$ws.Range("I132").Value = 2616.6086
$ws.Range("J132").Value = 3870
$ws.Range("K132").Value = 7849.825800000001
$ws.Range("L132").Value = 11610
$ws.Range("M132").Value = -5319.825800000001
$ws.Range("N132").Value = -16670
# Row 136
$ws.Range("H136").Value = 1641.8438
$ws.Range("I136").Value = 1504.473
$ws.Range("K136").Value = 4513.419
$ws.Range("M136").Value = -1963.419

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 3010.6667
$ws.Range("I122").Value = 828.5454999999999
$ws.Range("J122").Value = 3447.0908
$ws.Range("K122").Value = 7456.9095
$ws.Range("L122").Value = 31023.8172
$ws.Range("M122").Value = -5006.9095
$ws.Range("N122").Value = -35923.8172

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 12990.25
$ws.Range("J5").Value = 12990.25
$ws.Range("L5").Value = 12990.25
$ws.Range("N5").Value = -13214.25
# Row 107
$ws.Range("H107").Value = 4630229
$ws.Range("I107").Value = 311.58334
$ws.Range("J107").Value = 9260146
$ws.Range("K107").Value = 311.58334
$ws.Range("L107").Value = 9260146
$ws.Range("M107").Value = 1608.41666
$ws.Range("N107").Value = -9263986
# Row 132
$ws.Range("H132").Value = 4644.5293
$ws.Range("I132").Value = 3125.7144
$ws.Range("J132").Value = 5707.7
$ws.Range("K132").Value = 9377.143199999999
$ws.Range("L132").Value = 17123.1
$ws.Range("M132").Value = -6847.143199999999
$ws.Range("N132").Value = -22183.1

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 677.6
$ws.Range("I68").Value = 631.51544
$ws.Range("J68").Value = 2167.6667
$ws.Range("K68").Value = 631.51544
$ws.Range("L68").Value = 2167.6667
$ws.Range("M68").Value = 117.48456
$ws.Range("N68").Value = -3665.6667
# Row 71
$ws.Range("H71").Value = 677.6
$ws.Range("I71").Value = 631.51544
$ws.Range("J71").Value = 2167.6667
$ws.Range("K71").Value = 3157.5772
$ws.Range("L71").Value = 10838.3335
$ws.Range("M71").Value = 586.4227999999998
$ws.Range("N71").Value = -18326.3335
# Row 93
$ws.Range("H93").Value = 3473861.5
$ws.Range("I93").Value = 6537186
$ws.Range("J93").Value = 2093.5334
$ws.Range("K93").Value = 6537186
$ws.Range("L93").Value = 2093.5334
$ws.Range("M93").Value = -6535938
$ws.Range("N93").Value = -4589.5334

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 5131149
$ws.Range("I132").Value = 2866.7827
$ws.Range("J132").Value = 17546990
$ws.Range("K132").Value = 8600.348100000001
$ws.Range("L132").Value = 52640970
$ws.Range("M132").Value = -6070.348100000001
$ws.Range("N132").Value = -52646030
# Row 136
$ws.Range("H136").Value = 3624.4517
$ws.Range("I136").Value = 3911.4688
$ws.Range("J136").Value = 3318.3
$ws.Range("K136").Value = 11734.4064
$ws.Range("L136").Value = 9954.900000000001
$ws.Range("M136").Value = -9184.4064
$ws.Range("N136").Value = -15054.9
